$d = $word.ActiveDocument

$replacements = @(
    @("89×70=6230", "59×13=767"),
    @("62×84=5208", "16×99=1584"),
    @("16×91=1456", "70×72=5040"),
    @("41×20=820", "39×14=546"),
    @("27×40=1080", "21×63=1323"),
    @("99×55=5445", "60×30=1800"),
    @("96×32=3072", "67×40=2680"),
    @("77×35=2695", "63×52=3276"),
    @("96×46=4416", "70×24=1680"),
    @("90×47=4230", "28×66=1848"),
    @("14×61=854", "15×68=1020"),
    @("32×37=1184", "63×87=5481"),
    @("85×83=7055", "17×50=850"),
    @("65×35=2275", "16×31=496"),
    @("66×47=3102", "13×96=1248"),
    @("50×35=1750", "52×54=2808"),
    @("89×98=8722", "27×31=837"),
    @("60×33=1980", "12×32=384"),
    @("20×30=600", "69×91=6279"),
    @("86×39=3354", "32×13=416"),
    @("53×82=4346", "56×47=2632"),
    @("73×56=4088", "89×73=6497"),
    @("15×79=1185", "56×42=2352"),
    @("63×24=1512", "91×69=6279"),
    @("58×77=4466", "38×31=1178")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
